$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Add_Devices_LoopA")
$ws1.Columns.Item(5).ColumnWidth = 26.33203125
Write-Host "done"
